# Update the bills list: replace item/price data with a new, shorter list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (9 items instead of the previous 16)
$items = @(
    "45612 MILK WHOLE ",
    "22026 STIR FRY LARGE",
    "810873 FAMILY PACK TOMATO",
    "86247 RICE LG ",
    "727495 TEA GREEN ",
    "6023 COFFEE DECAF FD",
    "836067 PISTACHIOS SALTED",
    "42929 BREAD WHOLEMEAL",
    "Tortilla"
)

$prices = @(
    "1.55",
    "2.58",
    "1.39",
    "1.04",
    "0.65",
    "1.19",
    "1.99",
    "0.75",
    "0.99"
)

# Remove the rows that are no longer needed (old rows 11-17)
$ws.Range("A11:B17").Delete()

# Write the new item names and prices into rows 2-10
for ($i = 0; $i -lt $items.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $items[$i]
    $ws.Cells.Item($row, 2).Value = $prices[$i]
}
